# Updates the Price (D) and Volume(1h) (E) columns with a fresh snapshot,
# and fixes two rows whose Coin/Link/Price/Volume had been swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.796.37'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '3.851.39'
$ws.Range("E3").Value = '  -1.94%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = "'597.77"
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("D6").Value = "'167.78"
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '3.850.77'
$ws.Range("E7").Value = '  -1.85%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = '  -5.13%  '
$ws.Range("D11").Value = "'6.39"
$ws.Range("E11").Value = '  -0.80%  '
$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = '  -2.52%  '
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("D14").Value = "'36.71"
$ws.Range("E14").Value = '  -2.65%  '
$ws.Range("D15").Value = '4.508.50'
$ws.Range("E15").Value = '  -1.75%  '
$ws.Range("D16").Value = '3.865.24'
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("D17").Value = '68.016.02'
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").Value = "'18.06"
$ws.Range("E18").Value = '  +4.03%  '
$ws.Range("D19").Value = "'7.30"
$ws.Range("E19").Value = '  -2.60%  '
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").Value = "'10.74"
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").Value = "'463.81"
$ws.Range("E22").Value = '  -6.57%  '
$ws.Range("D23").Value = "'0.732"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -6.04%  '
$ws.Range("D25").Value = "'82.86"
$ws.Range("E25").Value = '  -2.75%  '
$ws.Range("D26").Value = "'2.22"
$ws.Range("E26").Value = '  -2.83%  '
$ws.Range("D27").Value = "'11.99"
$ws.Range("E27").Value = '  -1.52%  '
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").Value = "'9.92"
$ws.Range("E29").Value = '  -3.55%  '
$ws.Range("D30").Value = "'2.95"
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").Value = '4.013.31'
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("D32").Value = "'7.64"
$ws.Range("E32").Value = '  -1.98%  '
$ws.Range("D33").Value = "'2.30"
$ws.Range("E33").Value = '  -4.08%  '
$ws.Range("D34").Value = "'30.98"
$ws.Range("E34").Value = '  -3.27%  '
$ws.Range("D35").Value = "'9.42"
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("D36").Value = '3.825.16'
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("E37").Value = '  -3.16%  '
$ws.Range("D38").Value = "'3.67"
$ws.Range("E38").Value = '  +11.23%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").Value = "'1.02"
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = "'0.140"
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").Value = "'5.87"
$ws.Range("E41").Value = '  -2.45%  '
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").Value = "'0.310"
$ws.Range("E43").Value = '  -4.01%  '
$ws.Range("D44").Value = "'422.89"
$ws.Range("E44").Value = '  -2.38%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = "'1.96"
$ws.Range("E45").Value = '  -2.40%  '
$ws.Range("B46").Value = 'FLOKI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D46").Value = "'0.000295"
$ws.Range("E46").Value = '  +4.70%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").Value = "'47.07"
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = "'8.56"
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").Value = "'26.36"
$ws.Range("E50").Value = '  +2.44%  '
$ws.Range("D51").Value = "'142.14"
$ws.Range("E51").Value = '  -0.54%  '
